# This workbook contains averaged-intensity results for different sampling
# schemes. Three new "Spiral" schemes were run and their results inserted
# into the table (right after the existing "Gaussian-Quadrature" row), which
# pushes the remaining rows (NoRotation-tilt60deg ... HexGrid-60degTilt5degRes)
# down by three rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GossF")

# Rows 17-19 are brand new rows at the bottom of the table. Copy the
# formatting (bold/centered/bordered style used by the other index cells in
# column A) from an existing formatted cell before we populate the values.
$ws.Range("A10").Copy($ws.Range("A17"))
$ws.Range("A10").Copy($ws.Range("A18"))
$ws.Range("A10").Copy($ws.Range("A19"))

# Row 10: A=8, label="Gaussian-Quadrature"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9990556871472746
$ws.Range("D10").Value = 0.9568473323611687
$ws.Range("E10").Value = 1.009811656777114
$ws.Range("F10").Value = 0.9990556871472746
$ws.Range("G10").Value = 0.9686704768203884
$ws.Range("H10").Value = 1.038396101144211
$ws.Range("I10").Value = 1.008068444588058
$ws.Range("J10").Value = 0.9568473323611687
$ws.Range("K10").Value = 0.9833294945691415
$ws.Range("L10").Value = 0.991192590858208
$ws.Range("M10").Value = 0.9968082831397026

# Row 11: A=9, label="Spiral-90deg-10rot-5space"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9307027969858641
$ws.Range("D11").Value = 1.303916321476776
$ws.Range("E11").Value = 0.913226274997322
$ws.Range("F11").Value = 0.9307027969858641
$ws.Range("G11").Value = 1.16532510302329
$ws.Range("H11").Value = 0.7736688538060419
$ws.Range("I11").Value = 0.9164217406461607
$ws.Range("J11").Value = 1.303916321476776
$ws.Range("K11").Value = 1.108571298237049
$ws.Range("L11").Value = 1.019637047611456
$ws.Range("M11").Value = 1.000543515155909

# Row 12: A=10, label="Spiral-90deg-15rot-5space"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9302468709536855
$ws.Range("D12").Value = 1.305311981312687
$ws.Range("E12").Value = 0.912915563974888
$ws.Range("F12").Value = 0.9302468709536855
$ws.Range("G12").Value = 1.166061317612358
$ws.Range("H12").Value = 0.7729313415068712
$ws.Range("I12").Value = 0.916089629958706
$ws.Range("J12").Value = 1.305311981312687
$ws.Range("K12").Value = 1.109113772643788
$ws.Range("L12").Value = 1.019680321798737
$ws.Range("M12").Value = 1.000592784219866

# Row 13: A=11, label="Spiral-90deg-10rot-3space"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9306547013103774
$ws.Range("D13").Value = 1.304129990706679
$ws.Range("E13").Value = 0.9131621906631909
$ws.Range("F13").Value = 0.9306547013103774
$ws.Range("G13").Value = 1.165456013355459
$ws.Range("H13").Value = 0.7735119461771636
$ws.Range("I13").Value = 0.9163489234710389
$ws.Range("J13").Value = 1.304129990706679
$ws.Range("K13").Value = 1.108646090684935
$ws.Range("L13").Value = 1.019650395997656
$ws.Range("M13").Value = 1.000543960947318

# Row 14: A=12, label="NoRotation-tilt60deg"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.007764000000001
$ws.Range("D14").Value = 1.070868
$ws.Range("E14").Value = 0.9579519999999989
$ws.Range("F14").Value = 1.007764000000001
$ws.Range("G14").Value = 1.044072
$ws.Range("H14").Value = 0.8900959999999988
$ws.Range("I14").Value = 0.9731719999999981
$ws.Range("J14").Value = 1.070868
$ws.Range("K14").Value = 1.01441
$ws.Range("L14").Value = 1.011087
$ws.Range("M14").Value = 0.9906539999999994

# Row 15: A=13, label="Rotation-NoTilt"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.11
$ws.Range("D15").Value = 0.76
$ws.Range("E15").Value = 1.023325000000002
$ws.Range("F15").Value = 1.11
$ws.Range("G15").Value = 0.88
$ws.Range("H15").Value = 1.05
$ws.Range("I15").Value = 1.05
$ws.Range("J15").Value = 0.76
$ws.Range("K15").Value = 0.8916625000000009
$ws.Range("L15").Value = 1.000831250000001
$ws.Range("M15").Value = 0.9788875000000002

# Row 16: A=14, label="Rotation-60detTilt"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.062348833792
$ws.Range("D16").Value = 0.858053076377601
$ws.Range("E16").Value = 1.011931220582399
$ws.Range("F16").Value = 1.062348833792
$ws.Range("G16").Value = 0.9286418378751998
$ws.Range("H16").Value = 1.0278922395648
$ws.Range("I16").Value = 1.026796205670399
$ws.Range("J16").Value = 0.858053076377601
$ws.Range("K16").Value = 0.9349921484799999
$ws.Range("L16").Value = 0.998670491136
$ws.Range("M16").Value = 0.9859439023103999

# Row 17: A=15, label="HexGrid-90degTilt5degRes"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9954954706970016
$ws.Range("D17").Value = 0.9953554290253483
$ws.Range("E17").Value = 0.9954307713243284
$ws.Range("F17").Value = 0.9954954706970016
$ws.Range("G17").Value = 0.9952133550370761
$ws.Range("H17").Value = 0.9954841681819239
$ws.Range("I17").Value = 0.9955517101109677
$ws.Range("J17").Value = 0.9953554290253483
$ws.Range("K17").Value = 0.9953931001748384
$ws.Range("L17").Value = 0.99544428543592
$ws.Range("M17").Value = 0.9954218173961076

# Row 18: A=16, label="HexGrid-90degTilt22p5degRes"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9930793142625063
$ws.Range("D18").Value = 1.003742022620346
$ws.Range("E18").Value = 0.9937771487372318
$ws.Range("F18").Value = 0.9930793142625063
$ws.Range("G18").Value = 1.000030268885514
$ws.Range("H18").Value = 0.9928757061572363
$ws.Range("I18").Value = 0.9930979672522963
$ws.Range("J18").Value = 1.003742022620346
$ws.Range("K18").Value = 0.998759585678789
$ws.Range("L18").Value = 0.9959194499706476
$ws.Range("M18").Value = 0.9961004046525219

# Row 19: A=17, label="HexGrid-60degTilt5degRes"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9899634186265653
$ws.Range("D19").Value = 1.017094200222554
$ws.Range("E19").Value = 0.9906647134994233
$ws.Range("F19").Value = 0.9899634186265653
$ws.Range("G19").Value = 1.00794593767014
$ws.Range("H19").Value = 0.9821302494082897
$ws.Range("I19").Value = 0.9897450305347656
$ws.Range("J19").Value = 1.017094200222554
$ws.Range("K19").Value = 1.003879456860989
$ws.Range("L19").Value = 0.996921437743777
$ws.Range("M19").Value = 0.9962572583269563
